{"js": "// The \"ACTIVITIES\" table's second row (Past roles) has a date cell whose\n// last paragraph reads \"<27 spaces>August 2021 \u2013 Present\". The edit:\n//   1) shortens the leading whitespace run from 27 to 21 spaces, and\n//   2) changes the date \"August 2021\" to \"December 2022\" (keeping \" \u2013 Present\").\n// Final paragraph text becomes \"<21 spaces>December 2022 \u2013 Present\".\n\n// Step 1: shrink the whitespace run that immediately precedes \"August 2021\".\n// The combined string \"<27 spaces>August 2021\" is unique in the document, so\n// it safely anchors the edit without touching anything else.\nconst whitespacePrefix = \"                           \"; // 27 spaces\nconst whitespaceNew = \"                     \"; // 21 spaces\n\nconst prefixMatches = context.document.body.search(whitespacePrefix + \"August 2021\", { matchCase: true });\nprefixMatches.load(\"items\");\nawait context.sync();\n\nif (prefixMatches.items.length > 0) {\n  const prefixRange = prefixMatches.items[0];\n  prefixRange.insertText(whitespaceNew + \"August 2021\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// Step 2: change the date itself from \"August 2021\" to \"December 2022\",\n// leaving the trailing \" \u2013 Present\" untouched.\nconst dateMatches = context.document.body.search(\"August 2021\", { matchCase: true });\ndateMatches.load(\"items\");\nawait context.sync();\n\nif (dateMatches.items.length > 0) {\n  const dateRange = dateMatches.items[0];\n  dateRange.insertText(\"December 2022\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# The \"ACTIVITIES\" table's second row (Past roles) has a date cell whose\n# last paragraph reads \"<27 spaces>August 2021 - Present\" (en dash). The edit:\n#   1) shortens the leading whitespace run from 27 to 21 spaces, and\n#   2) changes the date \"August 2021\" to \"December 2022\" (keeping \" - Present\").\n# Final paragraph text becomes \"<21 spaces>December 2022 - Present\".\n\n$d = $word.ActiveDocument\n\n# Step 1: shrink the whitespace run that immediately precedes \"August 2021\"\n# from 27 spaces to 21 spaces. The combined string \"<27 spaces>August 2021\"\n# is unique in the document, so it safely anchors the edit.\n$whitespaceOld = \"                           \"\n$whitespaceNew = \"                     \"\n\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Text = $whitespaceOld + \"August 2021\"\n$find1.Replacement.ClearFormatting()\n$find1.Replacement.Text = $whitespaceNew + \"August 2021\"\n$find1.Execute([ref]$find1.Text, $true, $false, $false, $false, $false, $true, 1, $false, [ref]$find1.Replacement.Text, 2) | Out-Null\n\n# Step 2: change the date itself from \"August 2021\" to \"December 2022\",\n# leaving the trailing \" - Present\" untouched.\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Text = \"August 2021\"\n$find2.Replacement.ClearFormatting()\n$find2.Replacement.Text = \"December 2022\"\n$find2.Execute([ref]$find2.Text, $true, $false, $false, $false, $false, $true, 1, $false, [ref]$find2.Replacement.Text, 2) | Out-Null\n"}
